# Weekly CompStat refresh: new report week (Volume 30 No. 28, covering 7/10/2023-7/16/2023)
# and refreshed crime-complaint figures for the 69th Precinct.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: issue number and reporting week ---
$ws.Range("A8").Value2 = "Volume 30   Number  28"
$ws.Range("C9").Value2 = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# --- Crime-complaint figures (rows 14-30 = Murder..Hate Crimes, row 37 = historical Rape) ---
# Row 14
$ws.Range("A14").Value2 = "Murder"
$ws.Range("C14").Value2 = "0"
$ws.Range("D14").Value2 = 1
$ws.Range("D14").NumberFormat = '#,##0'
$ws.Range("E14").Value2 = -100
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F14").Value2 = "0"
$ws.Range("G14").Value2 = 1
$ws.Range("G14").NumberFormat = '#,##0'
$ws.Range("H14").Value2 = -100
$ws.Range("H14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I14").Value2 = 4
$ws.Range("J14").Value2 = 4
$ws.Range("K14").Value2 = 0
$ws.Range("L14").Value2 = -20
$ws.Range("M14").Value2 = 33.333333333333
$ws.Range("N14").Value2 = -33.333333333333

# Row 15
$ws.Range("A15").Value2 = "Rape"
$ws.Range("C15").Value2 = 2
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("D15").Value2 = "0"
$ws.Range("E15").Value2 = "***.*"
$ws.Range("F15").Value2 = 3
$ws.Range("G15").Value2 = 1
$ws.Range("H15").Value2 = 200
$ws.Range("I15").Value2 = 9
$ws.Range("J15").Value2 = 11
$ws.Range("K15").Value2 = -18.181818181818
$ws.Range("L15").Value2 = -10
$ws.Range("M15").Value2 = 0
$ws.Range("N15").Value2 = -35.714285714285

# Row 16
$ws.Range("A16").Value2 = "Robbery"
$ws.Range("C16").Value2 = 1
$ws.Range("D16").Value2 = 2
$ws.Range("E16").Value2 = -50
$ws.Range("F16").Value2 = 3
$ws.Range("G16").Value2 = 11
$ws.Range("H16").Value2 = -72.727272727272
$ws.Range("I16").Value2 = 69
$ws.Range("J16").Value2 = 65
$ws.Range("K16").Value2 = 6.153846153846
$ws.Range("L16").Value2 = 200
$ws.Range("M16").Value2 = -31.683168316831
$ws.Range("N16").Value2 = -79.705882352941

# Row 17
$ws.Range("A17").Value2 = "Fel. Assault"
$ws.Range("C17").Value2 = 5
$ws.Range("D17").Value2 = 3
$ws.Range("E17").Value2 = 66.666666666666
$ws.Range("F17").Value2 = 21
$ws.Range("G17").Value2 = 17
$ws.Range("H17").Value2 = 23.529411764705
$ws.Range("I17").Value2 = 121
$ws.Range("J17").Value2 = 130
$ws.Range("K17").Value2 = -6.923076923076
$ws.Range("L17").Value2 = 23.469387755102
$ws.Range("M17").Value2 = 42.35294117647
$ws.Range("N17").Value2 = -3.968253968253

# Row 18
$ws.Range("A18").Value2 = "Burglary"
$ws.Range("C18").Value2 = 1
$ws.Range("D18").Value2 = 1
$ws.Range("E18").Value2 = 0
$ws.Range("F18").Value2 = 2
$ws.Range("G18").Value2 = 6
$ws.Range("H18").Value2 = -66.666666666666
$ws.Range("I18").Value2 = 30
$ws.Range("J18").Value2 = 52
$ws.Range("K18").Value2 = -42.307692307692
$ws.Range("L18").Value2 = -23.076923076923
$ws.Range("M18").Value2 = -77.941176470588
$ws.Range("N18").Value2 = -90.963855421686

# Row 19
$ws.Range("A19").Value2 = "Gr. Larceny"
$ws.Range("C19").Value2 = 3
$ws.Range("D19").Value2 = 3
$ws.Range("E19").Value2 = 0
$ws.Range("F19").Value2 = 15
$ws.Range("G19").Value2 = 26
$ws.Range("H19").Value2 = -42.307692307692
$ws.Range("I19").Value2 = 114
$ws.Range("J19").Value2 = 135
$ws.Range("K19").Value2 = -15.555555555555
$ws.Range("L19").Value2 = 75.384615384615
$ws.Range("M19").Value2 = -1.724137931034
$ws.Range("N19").Value2 = -45.454545454545

# Row 20
$ws.Range("A20").Value2 = "G.L.A."
$ws.Range("C20").Value2 = 2
$ws.Range("D20").Value2 = 4
$ws.Range("D20").NumberFormat = '#,##0'
$ws.Range("E20").Value2 = -50
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F20").Value2 = 15
$ws.Range("G20").Value2 = 11
$ws.Range("H20").Value2 = 36.363636363636
$ws.Range("I20").Value2 = 85
$ws.Range("J20").Value2 = 81
$ws.Range("K20").Value2 = 4.938271604938
$ws.Range("L20").Value2 = 70
$ws.Range("M20").Value2 = 37.096774193548
$ws.Range("N20").Value2 = -92.517605633802

# Row 21
$ws.Range("A21").Value2 = "TOTAL"
$ws.Range("C21").Value2 = 14
$ws.Range("D21").Value2 = 14
$ws.Range("E21").Value2 = 0
$ws.Range("F21").Value2 = 59
$ws.Range("G21").Value2 = 73
$ws.Range("H21").Value2 = -19.17808219178
$ws.Range("I21").Value2 = 432
$ws.Range("J21").Value2 = 478
$ws.Range("K21").Value2 = -9.623430962343
$ws.Range("L21").Value2 = 48.965517241379
$ws.Range("M21").Value2 = -15.625
$ws.Range("N21").Value2 = -80.02773925104

# Row 22
$ws.Range("A22").Value2 = "Transit"
$ws.Range("C22").Value2 = "0"
$ws.Range("D22").Value2 = "0"
$ws.Range("E22").Value2 = "***.*"
$ws.Range("F22").Value2 = 3
$ws.Range("G22").Value2 = 1
$ws.Range("H22").Value2 = 200
$ws.Range("I22").Value2 = 6
$ws.Range("J22").Value2 = 6
$ws.Range("K22").Value2 = 0
$ws.Range("L22").Value2 = 0
$ws.Range("M22").Value2 = 0
$ws.Range("N22").Value2 = "***.*"

# Row 23
$ws.Range("A23").Value2 = "Housing"
$ws.Range("C23").Value2 = 2
$ws.Range("D23").Value2 = 3
$ws.Range("E23").Value2 = -33.333333333333
$ws.Range("F23").Value2 = 9
$ws.Range("G23").Value2 = 8
$ws.Range("H23").Value2 = 12.5
$ws.Range("I23").Value2 = 67
$ws.Range("J23").Value2 = 54
$ws.Range("K23").Value2 = 24.074074074074
$ws.Range("L23").Value2 = 26.415094339622
$ws.Range("M23").Value2 = 168
$ws.Range("N23").Value2 = "***.*"

# Row 24
$ws.Range("A24").Value2 = "Petit Larceny"
$ws.Range("C24").Value2 = 12
$ws.Range("D24").Value2 = 10
$ws.Range("E24").Value2 = 20
$ws.Range("F24").Value2 = 33
$ws.Range("G24").Value2 = 41
$ws.Range("H24").Value2 = -19.512195121951
$ws.Range("I24").Value2 = 317
$ws.Range("J24").Value2 = 288
$ws.Range("K24").Value2 = 10.069444444444
$ws.Range("L24").Value2 = 51.674641148325
$ws.Range("M24").Value2 = 11.228070175438
$ws.Range("N24").Value2 = "***.*"

# Row 25
$ws.Range("A25").Value2 = "Misd. Assault"
$ws.Range("C25").Value2 = 8
$ws.Range("D25").Value2 = 7
$ws.Range("E25").Value2 = 14.285714285714
$ws.Range("F25").Value2 = 34
$ws.Range("G25").Value2 = 33
$ws.Range("H25").Value2 = 3.030303030303
$ws.Range("I25").Value2 = 204
$ws.Range("J25").Value2 = 211
$ws.Range("K25").Value2 = -3.317535545023
$ws.Range("L25").Value2 = 34.210526315789
$ws.Range("M25").Value2 = -22.727272727272
$ws.Range("N25").Value2 = "***.*"

# Row 26
$ws.Range("A26").Value2 = "UCR Rape*"
$ws.Range("C26").Value2 = 2
$ws.Range("C26").NumberFormat = '#,##0'
$ws.Range("D26").Value2 = "0"
$ws.Range("E26").Value2 = "***.*"
$ws.Range("F26").Value2 = 4
$ws.Range("G26").Value2 = 2
$ws.Range("H26").Value2 = 100
$ws.Range("I26").Value2 = 12
$ws.Range("J26").Value2 = 18
$ws.Range("K26").Value2 = -33.333333333333
$ws.Range("L26").Value2 = -14.285714285714
$ws.Range("M26").Value2 = "***.*"
$ws.Range("N26").Value2 = "***.*"

# Row 27
$ws.Range("A27").Value2 = "Other Sex Crimes"
$ws.Range("C27").Value2 = 1
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("D27").Value2 = 1
$ws.Range("E27").Value2 = 0
$ws.Range("F27").Value2 = 1
$ws.Range("F27").NumberFormat = '#,##0'
$ws.Range("G27").Value2 = 3
$ws.Range("H27").Value2 = -66.666666666666
$ws.Range("I27").Value2 = 12
$ws.Range("J27").Value2 = 16
$ws.Range("K27").Value2 = -25
$ws.Range("L27").Value2 = -20
$ws.Range("M27").Value2 = "***.*"
$ws.Range("N27").Value2 = "***.*"

# Row 28
$ws.Range("A28").Value2 = "Shooting Vic."
$ws.Range("C28").Value2 = 2
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("D28").Value2 = 2
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value2 = 0
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F28").Value2 = 2
$ws.Range("F28").NumberFormat = '#,##0'
$ws.Range("G28").Value2 = 3
$ws.Range("H28").Value2 = -33.333333333333
$ws.Range("I28").Value2 = 11
$ws.Range("J28").Value2 = 18
$ws.Range("K28").Value2 = -38.888888888888
$ws.Range("L28").Value2 = -38.888888888888
$ws.Range("M28").Value2 = -15.384615384615
$ws.Range("N28").Value2 = -31.25

# Row 29
$ws.Range("A29").Value2 = "Shooting Inc."
$ws.Range("C29").Value2 = 2
$ws.Range("C29").NumberFormat = '#,##0'
$ws.Range("D29").Value2 = 2
$ws.Range("D29").NumberFormat = '#,##0'
$ws.Range("E29").Value2 = 0
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F29").Value2 = 2
$ws.Range("F29").NumberFormat = '#,##0'
$ws.Range("G29").Value2 = 3
$ws.Range("H29").Value2 = -33.333333333333
$ws.Range("I29").Value2 = 11
$ws.Range("J29").Value2 = 15
$ws.Range("K29").Value2 = -26.666666666666
$ws.Range("L29").Value2 = -21.428571428571
$ws.Range("M29").Value2 = 10
$ws.Range("N29").Value2 = -26.666666666666

# Row 30
$ws.Range("A30").Value2 = "Hate Crimes"
$ws.Range("C30").Value2 = "0"
$ws.Range("D30").Value2 = "0"
$ws.Range("E30").Value2 = "***.*"
$ws.Range("F30").Value2 = "0"
$ws.Range("G30").Value2 = "0"
$ws.Range("H30").Value2 = "***.*"
$ws.Range("I30").Value2 = "0"
$ws.Range("J30").Value2 = 1
$ws.Range("K30").Value2 = -100
$ws.Range("L30").Value2 = "***.*"
$ws.Range("M30").Value2 = "***.*"
$ws.Range("N30").Value2 = "***.*"

# Row 37
$ws.Range("A37").Value2 = "Rape"

Write-Host "done"